$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Q4_19_20")

# Insert a new row at row 5, pushing the existing rows 5-7 (A13/Rail, F9/RPE,
# Columbia/Rail) down to rows 6-8. This also grows the sheet's used range
# from B2:L7 to B2:L8.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new project entry.
$ws.Range("B5").Value = "A11"
$ws.Range("C5").Value = "HSMRPG"
